$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: map of cell address -> new value
# Percentage cells need NumberFormat "@" (Text) set first so Excel
# does not auto-convert the "NN%" string into a numeric percentage.
$percentCells = @("H2", "H5", "H9", "H19", "H20", "H23", "H24", "H25", "H28", "H30", "H34", "H37", "H39", "H41", "H45")
foreach ($pc in $percentCells) {
    $ws.Range($pc).NumberFormat = "@"
}

$updates = @{
    'E2' = '2026-02-18 21:48:26'
    'H2' = '68%'
    'I2' = '0.9 mm'
    'E3' = '2026-02-18 21:48:29'
    'I3' = '0.6 mm'
    'E4' = '2026-02-18 21:48:31'
    'J4' = '1012.6 hPa'
    'E5' = '2026-02-18 21:48:34'
    'H5' = '67%'
    'I5' = '1.2 mm'
    'O5' = '0.6 °C'
    'E6' = '2026-02-18 21:48:36'
    'J6' = '1012.3 hPa'
    'E7' = '2026-02-18 21:48:39'
    'J7' = '1013.9 hPa'
    'E8' = '2026-02-18 21:48:41'
    'J8' = '1013.6 hPa'
    'E9' = '2026-02-18 21:48:44'
    'H9' = '77%'
    'E10' = '2026-02-18 21:48:46'
    'E11' = '2026-02-18 21:48:49'
    'E12' = '2026-02-18 21:48:51'
    'O12' = '11.0 °C'
    'E13' = '2026-02-18 21:48:53'
    'J13' = '1014.9 hPa'
    'O13' = '4.0 °C'
    'E14' = '2026-02-18 21:48:56'
    'O14' = '12.3 °C'
    'E15' = '2026-02-18 21:48:58'
    'E16' = '2026-02-18 21:49:00'
    'G16' = '72 cm'
    'I16' = '2.3 mm'
    'O16' = '-0.2 °C'
    'E17' = '2026-02-18 21:49:03'
    'E18' = '2026-02-18 21:49:05'
    'J18' = '1012.8 hPa'
    'E19' = '2026-02-18 21:49:08'
    'H19' = '87%'
    'E20' = '2026-02-18 21:49:10'
    'H20' = '76%'
    'I20' = '0.2 mm'
    'O20' = '-0.5 °C'
    'E21' = '2026-02-18 21:49:13'
    'J21' = '1014.4 hPa'
    'E22' = '2026-02-18 21:49:15'
    'I22' = '1.3 mm'
    'N22' = '-3.9 °C 21:28 TU'
    'E23' = '2026-02-18 21:49:18'
    'H23' = '58%'
    'I23' = '0.4 mm'
    'N23' = '-3.0 °C 21:20 TU'
    'E24' = '2026-02-18 21:49:20'
    'H24' = '85%'
    'J24' = '1014.5 hPa'
    'L24' = '32.4 km/h - 290º 21:26 TU'
    'E25' = '2026-02-18 21:49:23'
    'H25' = '49%'
    'N25' = '-2.3 °C 21:14 TU'
    'O25' = '1.8 °C'
    'E26' = '2026-02-18 21:49:25'
    'J26' = '1011.8 hPa'
    'O26' = '5.5 °C'
    'E27' = '2026-02-18 21:49:28'
    'E28' = '2026-02-18 21:49:30'
    'H28' = '72%'
    'J28' = '1012.5 hPa'
    'E29' = '2026-02-18 21:49:33'
    'E30' = '2026-02-18 21:49:35'
    'H30' = '75%'
    'J30' = '1012.0 hPa'
    'O30' = '11.2 °C'
    'E31' = '2026-02-18 21:49:38'
    'J31' = '1010.8 hPa'
    'E32' = '2026-02-18 21:49:40'
    'E33' = '2026-02-18 21:49:43'
    'J33' = '1013.7 hPa'
    'E34' = '2026-02-18 21:49:45'
    'H34' = '48%'
    'O34' = '3.0 °C'
    'E35' = '2026-02-18 21:49:47'
    'J35' = '1014.1 hPa'
    'O35' = '9.2 °C'
    'E36' = '2026-02-18 21:49:50'
    'J36' = '1012.5 hPa'
    'O36' = '11.9 °C'
    'E37' = '2026-02-18 21:49:52'
    'H37' = '78%'
    'J37' = '1014.1 hPa'
    'E38' = '2026-02-18 21:49:55'
    'E39' = '2026-02-18 21:49:57'
    'H39' = '42%'
    'O39' = '1.2 °C'
    'E40' = '2026-02-18 21:50:00'
    'I40' = '1.6 mm'
    'J40' = '1015.0 hPa'
    'O40' = '6.4 °C'
    'E41' = '2026-02-18 21:50:02'
    'H41' = '87%'
    'J41' = '1014.1 hPa'
    'E42' = '2026-02-18 21:50:05'
    'E43' = '2026-02-18 21:50:07'
    'E44' = '2026-02-18 21:50:09'
    'E45' = '2026-02-18 21:50:12'
    'H45' = '65%'
    'I45' = '1.2 mm'
    'J45' = '1011.5 hPa'
    'E46' = '2026-02-18 21:50:15'
    'J46' = '1014.5 hPa'
    'K46' = '12.2 MJ/m2'
    'O46' = '10.9 °C'
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
